$d = $word.ActiveDocument

# The table row for the Gurage "e" series (row 16 of the first table) has cells:
#   1: "e"        row label
#   2: "<U+1260>ee"   base glyph + line break + "ee"
#   3: "<U+12A1>u"
#   4: "<U+12A2>i"
#   5: "<U+12A3>a"
#   6: "<U+12A4>ie"
#   7: "<U+1265>e"    base glyph + line break + "e"
#   8: "<U+12A6>o"
#
# Fix: cell 2's base glyph should be U+12A7 (ኧ), not U+1260 (አ), and the cell
# should only have a single trailing "e" (not "ee"). The extra "e" moves to
# cell 7, which becomes "<U+1265>ee".

$t = $d.Tables.Item(1)
$row = $t.Rows.Item(16)

# --- Cell 2: አee -> ኧe -------------------------------------------------
$cell2 = $row.Cells.Item(2)
$chars2 = $cell2.Range.Characters

# Character 1 is the base-glyph run (rFonts "Gurage Fider 1998" / cs
# "Abyssinica SIL test"). Replace its text in place so the run and its
# formatting are preserved, just with the corrected character.
$chars2.Item(1).Text = [string][char]0x12A7

# Character 4 is the second, duplicate "e" run; delete it so only the first
# "e" (character 3) remains.
$chars2.Item(4).Delete()

# --- Cell 7: እe -> እee --------------------------------------------------
$cell7 = $row.Cells.Item(7)
$chars7 = $cell7.Range.Characters

# Character 4 is the paragraph mark that ends the cell's text. Insert the
# extra "e" immediately before it, i.e. right after the existing "e".
$chars7.Item(4).InsertBefore("e")

$d.Save()
